$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra ticker row (MSFT) - table now only has AAPL
$ws.Rows(3).EntireRow.Delete()

# Remove the now-unused 3rd metric-group column (Current debt_3) - table now only
# has 4 metric groups x 2 columns instead of 3 metric groups x 3 columns
$ws.Columns("K").EntireColumn.Delete()

# Replace the header row with the corrected metric names picked up from the
# openbb API (revenue / free cash flow across TTM + 3 prior fiscal years)
$ws.Range("C1").Value = "Total revenue_TTM"
$ws.Range("D1").Value = "Free cash flow_TTM"
$ws.Range("E1").Value = "Total revenue_2022"
$ws.Range("F1").Value = "Free cash flow_2022"
$ws.Range("G1").Value = "Total revenue_2021"
$ws.Range("H1").Value = "Free cash flow_2021"
$ws.Range("I1").Value = "Total revenue_2020"
$ws.Range("J1").Value = "Free cash flow_2020"

# Replace the data row with the correct values for AAPL (also fills in the
# previously-missing EBITDA/placeholder cells)
$ws.Range("C2").Value = 383933000000
$ws.Range("D2").Value = 100987000000
$ws.Range("E2").Value = 394328000000
$ws.Range("F2").Value = 111443000000
$ws.Range("G2").Value = 365817000000
$ws.Range("H2").Value = 92953000000
$ws.Range("I2").Value = 274515000000
$ws.Range("J2").Value = 73365000000
